$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F14").Value = 6
$ws.Range("F15").Value = -2
$ws.Range("F18").Value = -3
$ws.Range("F21").Value = -4
$ws.Range("F22").Value = 1
$ws.Range("F23").Value = 1
$ws.Range("F30").Value = 5
$ws.Range("F32").Value = -4
$ws.Range("F34").Value = 2
$ws.Range("F35").Value = -4
$ws.Range("F39").Value = -4
$ws.Range("F42").Value = 4
$ws.Range("F43").Value = 6
$ws.Range("F46").Value = 2
$ws.Range("F47").Value = -1
$ws.Range("F48").Value = -2
$ws.Range("F50").Value = -3
$ws.Range("F51").Value = 1
$ws.Range("F52").Value = -1
$ws.Range("F54").Value = -3
$ws.Range("F55").Value = 1
$ws.Range("F56").Value = 4
$ws.Range("F57").Value = 3
